$d = $word.ActiveDocument

# Step 1: split the paragraph that ends with "that show PNGs of tiles. "
# into two paragraphs, both sharing the same ListParagraph/numId=7 list
# formatting. The new (second) paragraph starts out empty; the document's
# trailing _GoBack bookmark naturally relocates to the start of that new
# paragraph.
$d.Content.Find.Execute("that show PNGs of tiles. ", $true, $false, $false, `
    $false, $false, $true, 1, $false, `
    "that show PNGs of tiles. ^p", 2) | Out-Null

# Step 2: insert the new bullet text right before the (now relocated)
# _GoBack bookmark. Using InsertBefore on the bookmark's own collapsed
# range leaves the bookmark anchored after the inserted text, i.e. at the
# end of the new paragraph, matching the original "trailing bookmark"
# position relative to the paragraph's text.
$bm = $d.Bookmarks("_GoBack")
$insertionPoint = $d.Range($bm.Start, $bm.Start)
$insertionPoint.InsertBefore( `
    "Created custom ImageView class to drawText a number on top of a Sun Tile")
